# Updates the "cryptos" sheet with refreshed price / volume(1h) figures
# for Sun Jan  8 23:22:42 UTC 2023 (GitHub Actions symbol-list refresh),
# plus the BKEXToken / KickToken row re-order (rows 41-42 swapped places).
#
# Price (column D) and Volume 1h (column E) are stored as literal text in
# this workbook (e.g. "273.87", "4.80%"), so each assignment is prefixed
# with a leading apostrophe to force Excel to keep the value as text
# instead of auto-converting it to a number/percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'273.87"
$ws.Range("E2").Value = "'4.80%"

# Row 3
$ws.Range("D3").Value = "'26.83"
$ws.Range("E3").Value = "'-1.21%"

# Row 4
$ws.Range("D4").Value = "'4.731"
$ws.Range("E4").Value = "'0.51%"

# Row 5
$ws.Range("D5").Value = "'0.06132"
$ws.Range("E5").Value = "'-1.14%"

# Row 6
$ws.Range("D6").Value = "'6.736"
$ws.Range("E6").Value = "'0.09%"

# Row 7
$ws.Range("D7").Value = "'0.8615"
$ws.Range("E7").Value = "'1.21%"

# Row 8
$ws.Range("D8").Value = "'0.9098"
$ws.Range("E8").Value = "'0.34%"

# Row 9
$ws.Range("D9").Value = "'0.1443"
$ws.Range("E9").Value = "'2.84%"

# Row 10
$ws.Range("D10").Value = "'0.05080"
$ws.Range("E10").Value = "'6.67%"

# Row 11
$ws.Range("D11").Value = "'0.07148"
$ws.Range("E11").Value = "'0.72%"

# Row 12
$ws.Range("D12").Value = "'0.03187"
$ws.Range("E12").Value = "'0.39%"

# Row 13
$ws.Range("D13").Value = "'0.09037"
$ws.Range("E13").Value = "'-0.25%"

# Row 14
$ws.Range("D14").Value = "'0.001536"
$ws.Range("E14").Value = "'0.35%"

# Row 15
$ws.Range("D15").Value = "'0.0006092"
$ws.Range("E15").Value = "'-0.82%"

# Row 16
$ws.Range("D16").Value = "'0.005932"
$ws.Range("E16").Value = "'-3.45%"

# Row 17
$ws.Range("D17").Value = "'3.462"
$ws.Range("E17").Value = "'-0.18%"

# Row 18
$ws.Range("D18").Value = "'3.183"
$ws.Range("E18").Value = "'0.36%"

# Row 19
$ws.Range("D19").Value = "'2.264"
$ws.Range("E19").Value = "'4.00%"

# Row 20
$ws.Range("D20").Value = "'0.3085"
$ws.Range("E20").Value = "'-0.68%"

# Row 21
$ws.Range("D21").Value = "'0.1300"
$ws.Range("E21").Value = "'1.51%"

# Row 22
$ws.Range("D22").Value = "'3.828"
$ws.Range("E22").Value = "'-7.18%"

# Row 23
$ws.Range("D23").Value = "'0.04222"
$ws.Range("E23").Value = "'0.04%"

# Row 24
$ws.Range("D24").Value = "'0.001176"
$ws.Range("E24").Value = "'-3.53%"

# Row 25
$ws.Range("D25").Value = "'0.004192"
$ws.Range("E25").Value = "'1.79%"

# Row 26
$ws.Range("D26").Value = "'0.0001196"
$ws.Range("E26").Value = "'-0.42%"

# Row 27
$ws.Range("D27").Value = "'0.0001676"
$ws.Range("E27").Value = "'3.69%"

# Row 40
$ws.Range("D40").Value = "'0.03980"
$ws.Range("E40").Value = "'2.04%"

# Row 41
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006201"
$ws.Range("E41").Value = "'50.11%"

# Row 42
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1130"
$ws.Range("E42").Value = "'1.48%"

# Row 43
$ws.Range("D43").Value = "'0.002164"
$ws.Range("E43").Value = "'-0.90%"

# Row 44
$ws.Range("D44").Value = "'0.01195"
$ws.Range("E44").Value = "'-11.20%"

# Row 45
$ws.Range("D45").Value = "'0.00005148"
$ws.Range("E45").Value = "'-0.48%"

# Row 46
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'-0.11%"

# Row 47
$ws.Range("D47").Value = "'0.8985"
$ws.Range("E47").Value = "'429.22%"

# Row 48
$ws.Range("D48").Value = "'0.02991"
$ws.Range("E48").Value = "'-16.70%"

# Row 49
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("E49").Value = "'-0.11%"

# Row 50
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("E50").Value = "'-0.11%"
